$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coordenadaX (E) and coordenadaY (F) values for rows 2-5
# using linear interpolator and avg error (per commit message)
$ws.Range("E2").Value = 792.927
$ws.Range("F2").Value = -70.439016520271

$ws.Range("E3").Value = 383.065
$ws.Range("F3").Value = 1018.382359662934

$ws.Range("E4").Value = 268.239
$ws.Range("F4").Value = 16226.40414747018

$ws.Range("E5").Value = 642.051
$ws.Range("F5").Value = -82696.38293800216
